$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) and Volume (E) columns for the cells
# we touch, so numeric-looking strings (e.g. '1.00', '6.46') are not
# auto-coerced into Excel numbers and lose their original formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.574.91'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.357.38'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.36'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.07'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.46%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.356.17'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.33'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.119'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.940.56'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.99%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.56'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.732.81'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.364.92'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.85'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.94'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.63%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.992'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.89'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000112'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +14.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.32'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.174'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.95'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.86%  '

$ws.Range("B31").Value = 'RenderToken'

$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.46'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.12%  '

$ws.Range("B32").Value = 'PancakeSwap'

$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.96'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.95%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.60%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.85'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.64'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.51'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.57%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0758'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.884.86'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.69'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.75%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0312'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.83'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.78%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.738'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.70'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.04'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +13.36%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.819'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.85%  '
